$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 488; existing rows 488-538 shift down to 489-539.
$ws.Rows("488:488").Insert()

# Populate the newly inserted row 488 with the new data record.
$ws.Range("A488").Value = 4
$ws.Range("B488").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C488").Value = "Los Lagos"
$ws.Range("D488").Value = 45212
$ws.Range("E488").Value = 10
$ws.Range("F488").Value = 100112040
$ws.Range("G488").Value = "Cilantro"
$ws.Range("H488").Value = "Sin especificar"
$ws.Range("I488").Value = "Primera"
$ws.Range("J488").Value = 180
$ws.Range("K488").Value = 12000
$ws.Range("L488").Value = 12000
$ws.Range("M488").Value = 12000
$ws.Range("N488").Value = "$/caja 36 atados"
$ws.Range("O488").Value = "Región Metropolitana"
$ws.Range("P488").Value = 333
$ws.Range("Q488").Value = 36
$ws.Range("R488").Value = "Hortaliza"

# Make sure the D488 cell keeps the date number format used by the rest of column D.
$ws.Range("D488").NumberFormat = $ws.Range("D489").NumberFormat()
